# issue #5: stock data output to json file
#
# Adds a new "property_category" column (value "stock") to the 股票 (stock)
# worksheet, and fixes two company names that had a stray internal space.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new blank column before the existing "date" column (column H),
# shifting date/legislator_name/legislator_id one column to the right.
$ws.Columns.Item(8).Insert()

# New header + values for the inserted "property_category" column.
$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
$ws.Range("H4").Value = "stock"

# Fix company names that had an erroneous embedded space.
$ws.Range("B2").Value = "春源鋼鐵工業股份有限公司"
$ws.Range("B4").Value = "中華開發金融控股股份有限公司"
